$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 101, shifting existing rows 101-136 down to 102-137.
$ws.Rows(101).Insert()

# Populate the new row 101 with the new record.
$ws.Cells.Item(101, 1).Value = 7
$ws.Cells.Item(101, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(101, 3).Value = "Ñuble"
$ws.Cells.Item(101, 4).Value = 44572
$ws.Cells.Item(101, 5).Value = 16
$ws.Cells.Item(101, 6).Value = 100112028
$ws.Cells.Item(101, 7).Value = "Sandia"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 800
$ws.Cells.Item(101, 11).Value = 2000
$ws.Cells.Item(101, 12).Value = 2300
$ws.Cells.Item(101, 13).Value = 2150
$ws.Cells.Item(101, 14).Value = "$/unidad"
$ws.Cells.Item(101, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(101, 16).Value = 2150
$ws.Cells.Item(101, 17).Value = 1
$ws.Cells.Item(101, 18).Value = "Hortaliza"
